$d = $word.ActiveDocument

# 1. Garments category bullet: append ",quantity" after the material list.
$d.Content.Find.Execute(
    "size and material (gents,ladies,unisex)", $true, $false, $false, $false, $false,
    $true, 1, $false, "size and material (gents,ladies,unisex),quantity", 2) | Out-Null

# 2. Electronics category bullet: append ",quantity" after the wattage text.
$d.Content.Find.Execute(
    "price, size and warranty, electricity wattage", $true, $false, $false, $false, $false,
    $true, 1, $false, "price, size and warranty, electricity wattage,quantity", 2) | Out-Null

# 3. Add two new bullet items under "Coding Guidelines", right after the
#    "Upload the Project to GitHub..." list item. InsertParagraphAfter on
#    that item's range creates a new paragraph that inherits the same
#    ListParagraph style / numbering, matching the surrounding bullets.
$rngUpload = $d.Content
$rngUpload.Find.Execute("Upload the Project to GitHub Link that will be shared to you.") | Out-Null
$rngUpload.InsertParagraphAfter()

$paras = $d.Paragraphs
$uploadIndex = 0
for ($i = 1; $i -le $paras.Count; $i++) {
    if ($paras.Item($i).Range.Text.TrimEnd() -eq "Upload the Project to GitHub Link that will be shared to you.") {
        $uploadIndex = $i
        break
    }
}

$logPara = $paras.Item($uploadIndex + 1)
$logPara.Range.Text = "Log message to a file"

$logPara.Range.InsertParagraphAfter()

$paras2 = $d.Paragraphs
$testPara = $paras2.Item($uploadIndex + 2)
$testPara.Range.Text = "Write atheist Five Unit Test cases"
